# Automatische test-sync: 2025-06-29 14:32:50
# Adds a new log row (#13) to the "Logs" sheet for the testmail that came in
# at 2025-06-29 14:32:40, keeps the conditional-formatting ranges in sync
# with the newly grown data range, and bumps the "Openingstijden / Locatie"
# tally on the "Dashboard" sheet from 3 to 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Logs")
$dash = $wb.Sheets.Item("Dashboard")

$nl = "`n"

# --- New row 13 ------------------------------------------------------------
$ws.Range("A13").Value = "Wanneer zijn jullie open?"
$ws.Range("B13").Value = "mailmind.test@zohomail.eu"
$ws.Range("C13").Value = "Testmail #1: Wanneer zijn jullie open?"
$ws.Range("D13").Value = "Openingstijden / Locatie"
$ws.Range("E13").Value = (
    "Beste klant," + $nl +
    "Bedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. We zijn gesloten in het weekend. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen." + $nl +
    "Met vriendelijke groet," + $nl +
    "[Naam bedrijf]"
)
$ws.Range("F13").Value = "2025-06-29 14:32:40"
$ws.Range("G13").Value = "Ja"
$ws.Range("H13").Value = "Nee"
$ws.Range("I13").Value = "Ja"

# Entering an embedded line break auto-grows the row; put it back to the
# sheet's standard height so row 13 matches the other (non-customised) rows.
$ws.Rows.Item(13).AutoFit()

# --- Keep conditional formatting ranges in sync with the new row ----------
$ws.Range("D2:D12").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D13"))
$ws.Range("G2:G12").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G13"))
$ws.Range("H2:H12").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H13"))
$ws.Range("I2:I12").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I13"))

# --- Dashboard tally: "Openingstijden / Locatie" count 3 -> 4 --------------
$dash.Range("B2").Value = 4
